$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1280.8572
$ws.Range("I33").Value = 1229.2941
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 1229.2941
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -1000.2941
$ws.Range("N33").Value = -1958
$ws.Range("H43").Value = 210230.9
$ws.Range("I43").Value = 8666.666999999999
$ws.Range("K43").Value = 8666.666999999999
$ws.Range("M43").Value = -8597.666999999999
$ws.Range("H63").Value = 75474
$ws.Range("I63").Value = 35948
$ws.Range("J63").Value = 115000
$ws.Range("K63").Value = 35948
$ws.Range("L63").Value = 115000
$ws.Range("M63").Value = -35324
$ws.Range("N63").Value = -116248
$ws.Range("H66").Value = 75474
$ws.Range("I66").Value = 35948
$ws.Range("J66").Value = 115000
$ws.Range("K66").Value = 107844
$ws.Range("L66").Value = 345000
$ws.Range("M66").Value = -104724
$ws.Range("N66").Value = -351240
$ws.Range("H137").Value = 3710.182
$ws.Range("I137").Value = 2482
$ws.Range("K137").Value = 7446
$ws.Range("M137").Value = -4896
$ws.Range("H138").Value = 2852.4185
$ws.Range("J138").Value = 3227.4675
$ws.Range("L138").Value = 9682.4025
$ws.Range("N138").Value = -19962.4025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13161076
$ws.Range("I32").Value = 14707853
$ws.Range("J32").Value = 13474.5
$ws.Range("K32").Value = 14707853
$ws.Range("L32").Value = 13474.5
$ws.Range("M32").Value = -14707566
$ws.Range("N32").Value = -14048.5
$ws.Range("H45").Value = 3169.111
$ws.Range("I45").Value = 2931.7144
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 2931.7144
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -2554.7144
$ws.Range("N45").Value = -4754
$ws.Range("H74").Value = 9268355
$ws.Range("I74").Value = 13890883
$ws.Range("K74").Value = 13890883
$ws.Range("M74").Value = -13890009
$ws.Range("H77").Value = 9268355
$ws.Range("I77").Value = 13890883
$ws.Range("K77").Value = 69454415
$ws.Range("M77").Value = -69450047
$ws.Range("H97").Value = 2205.5
$ws.Range("I97").Value = 2657.3333
$ws.Range("J97").Value = 1527.75
$ws.Range("K97").Value = 2657.3333
$ws.Range("L97").Value = 1527.75
$ws.Range("M97").Value = -2161.3333
$ws.Range("N97").Value = -2519.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5873.5
$ws.Range("I20").Value = 5873.5
$ws.Range("K20").Value = 5873.5
$ws.Range("M20").Value = -5626.5
$ws.Range("H86").Value = 3105.5625
$ws.Range("I86").Value = 2940.2307
$ws.Range("K86").Value = 2940.2307
$ws.Range("M86").Value = -1817.2307
$ws.Range("H89").Value = 3105.5625
$ws.Range("I89").Value = 2940.2307
$ws.Range("K89").Value = 14701.1535
$ws.Range("M89").Value = -9085.1535
$ws.Range("H105").Value = 2233.6667
$ws.Range("I105").Value = 1845
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 1845
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -98
$ws.Range("N105").Value = -6505
$ws.Range("H122").Value = 63996.668
$ws.Range("J122").Value = 63996.668
$ws.Range("L122").Value = 63996.668
$ws.Range("N122").Value = -73796.66800000001
$ws.Range("H130").Value = 61374.5
$ws.Range("J130").Value = 61374.5
$ws.Range("L130").Value = 61374.5
$ws.Range("N130").Value = -71414.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 22834.834
$ws.Range("I10").Value = 6333.6665
$ws.Range("J10").Value = 39336
$ws.Range("K10").Value = 6333.6665
$ws.Range("L10").Value = 39336
$ws.Range("M10").Value = -6194.6665
$ws.Range("N10").Value = -39614
$ws.Range("H31").Value = 1026751
$ws.Range("I31").Value = 33804.4
$ws.Range("J31").Value = 1302569.5
$ws.Range("K31").Value = 33804.4
$ws.Range("L31").Value = 1302569.5
$ws.Range("M31").Value = -33509.4
$ws.Range("N31").Value = -1303159.5
$ws.Range("H34").Value = 1026751
$ws.Range("I34").Value = 33804.4
$ws.Range("J34").Value = 1302569.5
$ws.Range("K34").Value = 33804.4
$ws.Range("L34").Value = 1302569.5
$ws.Range("M34").Value = -33602.4
$ws.Range("N34").Value = -1302973.5
$ws.Range("H94").Value = 4174.25
$ws.Range("I94").Value = 3234.5
$ws.Range("J94").Value = 4738.1
$ws.Range("K94").Value = 3234.5
$ws.Range("L94").Value = 4738.1
$ws.Range("M94").Value = -2783.5
$ws.Range("N94").Value = -5640.1
$ws.Range("H99").Value = 3193.4375
$ws.Range("I99").Value = 3233.3333
$ws.Range("K99").Value = 3233.3333
$ws.Range("M99").Value = -1735.3333
$ws.Range("H108").Value = 95815.5
$ws.Range("J108").Value = 95815.5
$ws.Range("L108").Value = 95815.5
$ws.Range("N108").Value = -103495.5
$ws.Range("H122").Value = 2983.75
$ws.Range("I122").Value = 2983.75
$ws.Range("K122").Value = 8951.25
$ws.Range("M122").Value = -6501.25
$ws.Range("H126").Value = 3193.4375
$ws.Range("I126").Value = 3233.3333
$ws.Range("K126").Value = 9699.999899999999
$ws.Range("M126").Value = -7229.999899999999
$ws.Range("H127").Value = 52000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 52000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 52000
$ws.Range("N127").Value = -61920
$ws.Range("M127").ClearContents()
$ws.Range("H132").Value = 4833.3335
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 28503.25
$ws.Range("I51").Value = 14998
$ws.Range("K51").Value = 44994
$ws.Range("M51").Value = -44534
$ws.Range("H80").Value = 4933.222
$ws.Range("I80").Value = 3449.5
$ws.Range("J80").Value = 5357.143
$ws.Range("K80").Value = 10348.5
$ws.Range("L80").Value = 16071.429
$ws.Range("M80").Value = -9412.5
$ws.Range("N80").Value = -17943.429
$ws.Range("H83").Value = 4933.222
$ws.Range("I83").Value = 3449.5
$ws.Range("J83").Value = 5357.143
$ws.Range("K83").Value = 31045.5
$ws.Range("L83").Value = 48214.287
$ws.Range("M83").Value = -26365.5
$ws.Range("N83").Value = -57574.287
$ws.Range("H119").Value = 11571.429
$ws.Range("I119").Value = 3000
$ws.Range("K119").Value = 9000
$ws.Range("M119").Value = -4162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 17365.666
$ws.Range("J63").Value = 17365.666
$ws.Range("L63").Value = 17365.666
$ws.Range("N63").Value = -18737.666
$ws.Range("H66").Value = 17365.666
$ws.Range("J66").Value = 17365.666
$ws.Range("L66").Value = 52096.99800000001
$ws.Range("N66").Value = -58960.99800000001
$ws.Range("H80").Value = 4162.6665
$ws.Range("I80").Value = 3995.2
$ws.Range("K80").Value = 3995.2
$ws.Range("M80").Value = -2997.2
$ws.Range("H83").Value = 4162.6665
$ws.Range("I83").Value = 3995.2
$ws.Range("K83").Value = 19976
$ws.Range("M83").Value = -14984
$ws.Range("H132").Value = 333386660
$ws.Range("J132").Value = 80000
$ws.Range("L132").Value = 240000
$ws.Range("N132").Value = -245060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 20000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -19437
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -19853
$ws.Range("N49").ClearContents()
$ws.Range("H136").Value = 69347.37
$ws.Range("I136").Value = 12399.889
$ws.Range("J136").Value = 120600.1
$ws.Range("K136").Value = 37199.667
$ws.Range("L136").Value = 361800.3
$ws.Range("M136").Value = -34649.667
$ws.Range("N136").Value = -366900.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 945.7143
$ws.Range("I136").Value = 988.9231
$ws.Range("J136").Value = 384
$ws.Range("K136").Value = 2966.7693
$ws.Range("L136").Value = 1152
$ws.Range("M136").Value = -416.7692999999999
$ws.Range("N136").Value = -6252
